# Append the latest "GetConfigReq" snapshot row (row 30) to the sheet,
# mirroring the layout of the existing rows (e.g. row 29) with the new
# poll's timestamp/epoch and wifi RSSI reading.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "2022-03-21 17:09:30"
$ws.Range("B30").Value = ""
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = "NONE"
$ws.Range("E30").Value = "NONE"
$ws.Range("F30").Value = "CMNET"
$ws.Range("G30").Value = 223
$ws.Range("H30").Value = 5
$ws.Range("I30").Value = 5
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 1
$ws.Range("Q30").Value = "10F872226797"
$ws.Range("R30").Value = 0
$ws.Range("S30").Value = 0
$ws.Range("T30").Value = 0
$ws.Range("U30").Value = 0
$ws.Range("V30").Value = 0
$ws.Range("W30").Value = ""
$ws.Range("X30").Value = 0
$ws.Range("Y30").Value = 0
$ws.Range("Z30").Value = 0
$ws.Range("AA30").Value = 72
$ws.Range("AB30").Value = 77
$ws.Range("AC30").Value = 114
$ws.Range("AD30").Value = 34
$ws.Range("AE30").Value = 103
$ws.Range("AF30").Value = 151
$ws.Range("AG30").Value = "NONE"
$ws.Range("AH30").Value = "NONE"
$ws.Range("AI30").Value = 1
$ws.Range("AJ30").Value = 3600
$ws.Range("AK30").Value = "dataeu.hoymiles.com"
$ws.Range("AL30").Value = 1
$ws.Range("AM30").Value = 10081
$ws.Range("AN30").Value = 0
$ws.Range("AO30").Value = 0
$ws.Range("AP30").Value = 0
$ws.Range("AQ30").Value = 0
$ws.Range("AR30").Value = 0
$ws.Range("AS30").Value = 0
$ws.Range("AT30").Value = 1647878973
$ws.Range("AU30").Value = 0
$ws.Range("AV30").Value = 0
$ws.Range("AW30").Value = 0
$ws.Range("AX30").Value = 0
$ws.Range("AY30").Value = 0
$ws.Range("AZ30").Value = 0
$ws.Range("BA30").Value = 0
$ws.Range("BB30").Value = 0
$ws.Range("BC30").Value = 0
$ws.Range("BD30").Value = 0
$ws.Range("BE30").Value = "0negawsklov0negawsklov"
$ws.Range("BF30").Value = 58
$ws.Range("BG30").Value = "HomeSweetHome"
$ws.Range("BH30").Value = 0
$ws.Range("BI30").Value = 0
